$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: fill in the calculated percentage amounts for the base value in C5
# (756), one formula per percentage column defined in row 4 (D4:H4), using
# the text-percent multiplication style already used in this workbook.
$ws.Range("D5").Formula = '=756*"8%"'
$ws.Range("E5").Formula = '=756*"10%"'
$ws.Range("F5").Formula = '=756*"5%"'
$ws.Range("G5").Formula = '=756*"3%"'
$ws.Range("H5").Formula = '=756*"12%"'

# Student identification block (Numara / Ad Soyad / Bölüm).
$ws.Range("G14").Value = 20215070055
$ws.Range("G15").Value = "Muhammed Ali Harmancı"
$ws.Range("G16").Value = "Yönetim Bilişim Sistemleri"

# Restore the saved view/selection state (scrolled so column C is left-most,
# with the active cell on G17).
try {
    $excel.ActiveWindow.ScrollColumn = 3
    $excel.ActiveWindow.ScrollRow = 1
} catch {
}
$ws.Range("G17").Select()
